$d = $word.ActiveDocument

# --- Paragraph 1: "This is a work in progress!" -> "Python for Geospatial
#     Big Data and Data Science Using the FASRC", set in Courier New
#     (both the run and the paragraph mark get the new font). ---
$p1 = $d.Paragraphs(1)
$body1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$body1.Text = "Python for Geospatial Big Data and Data Science Using the FASRC"

$p1 = $d.Paragraphs(1)
$p1.Range.Font.Name = "Courier New"
$p1.Range.Font.NameBi = "Courier New"

# --- Paragraph 2 is the blank <w:p/> separator - left untouched. ---

# --- Paragraph 3: "Workshop: Python for Geospatial Big Data and Data
#     Science Using the FASRC" -> emptied out and restyled as the
#     document Title (the heading text now lives solely in paragraph 4,
#     "Command Cheat Sheet", which keeps the Title style already). ---
$p3 = $d.Paragraphs(3)
$body3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$body3.Text = ""

$p3 = $d.Paragraphs(3)
$p3.Style = "Title"

Write-Output "ok"
